$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.977.56'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.956.89'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.18'
$ws.Range("E5").Value = '  +13.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.61'
$ws.Range("E6").Value = '  +11.91%  '
$ws.Range("E7").Value = '  -2.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.749'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("E10").Value = '  +4.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '55.95'
$ws.Range("E11").Value = '  +4.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000333'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.07'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.587.82'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.955.31'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("E16").Value = '  +2.87%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.96'
$ws.Range("E17").Value = '  -2.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.37'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.790.62'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.131'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.98'
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("E22").Value = '  +14.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.22'
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("E24").Value = '  -4.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.12'
$ws.Range("E25").Value = '  -3.68%  '
$ws.Range("E26").Value = '  -7.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.06'
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.40'
$ws.Range("E29").Value = '  -3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.86'
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.92'
$ws.Range("E31").Value = '  -4.44%  '
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0000102'
$ws.Range("E34").Value = '  +15.09%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '47.75'
$ws.Range("E35").Value = '  -3.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '69.46'
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '637.60'
$ws.Range("E37").Value = '  -6.96%  '
$ws.Range("E38").Value = '  -5.94%  '
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("B40").Value = 'Dai'
$ws.Range("C40").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.146'
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.74'
$ws.Range("E43").Value = '  -3.79%  '
$ws.Range("E44").Value = '  -2.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.18'
$ws.Range("E45").Value = '  -5.71%  '
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.41'
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.86'
$ws.Range("E49").Value = '  +23.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.864.94'
$ws.Range("E50").Value = '  +4.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.00'
$ws.Range("E51").Value = '  -3.12%  '
